$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '30.204.65'
Set-TextValue $ws.Range('E2') '  +0.09%  '

Set-TextValue $ws.Range('D3') '1.858.31'
Set-TextValue $ws.Range('E3') '  +0.05%  '

Set-TextValue $ws.Range('E4') '  +0.05%  '

Set-TextValue $ws.Range('D5') '235.40'
Set-TextValue $ws.Range('E5') '  +0.57%  '

Set-TextValue $ws.Range('D7') '0.4666'
Set-TextValue $ws.Range('E7') '  -0.16%  '

Set-TextValue $ws.Range('D8') '0.2841'
Set-TextValue $ws.Range('E8') '  +0.77%  '

Set-TextValue $ws.Range('D9') '0.06506'
Set-TextValue $ws.Range('E9') '  -0.72%  '

Set-TextValue $ws.Range('D10') '21.43'
Set-TextValue $ws.Range('E10') '  +6.55%  '

Set-TextValue $ws.Range('D11') '0.07897'
Set-TextValue $ws.Range('E11') '  +0.84%  '

Set-TextValue $ws.Range('D12') '96.95'
Set-TextValue $ws.Range('E12') '  +0.21%  '

Set-TextValue $ws.Range('D13') '1.864.08'
Set-TextValue $ws.Range('E13') '  +0.30%  '

Set-TextValue $ws.Range('E14') '  +0.80%  '

Set-TextValue $ws.Range('D15') '0.6753'
Set-TextValue $ws.Range('E15') '  +1.36%  '

Set-TextValue $ws.Range('D16') '277.67'
Set-TextValue $ws.Range('E16') '  -2.08%  '

Set-TextValue $ws.Range('D17') '30.210.50'
Set-TextValue $ws.Range('E17') '  +0.02%  '

Set-TextValue $ws.Range('E18') '  +7.96%  '

Set-TextValue $ws.Range('D19') '1.001'
Set-TextValue $ws.Range('E19') '  +0.05%  '

Set-TextValue $ws.Range('D20') '5.362'
Set-TextValue $ws.Range('E20') '  -1.84%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D21') '0.000007294'
Set-TextValue $ws.Range('E21') '  +0.70%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D22') '2.101.61'
Set-TextValue $ws.Range('E22') '  -0.22%  '

Set-TextValue $ws.Range('D23') '1.002'
Set-TextValue $ws.Range('E23') '  +0.08%  '

Set-TextValue $ws.Range('D24') '6.121'
Set-TextValue $ws.Range('E24') '  -0.37%  '

Set-TextValue $ws.Range('D25') '166.47'
Set-TextValue $ws.Range('E25') '  -0.65%  '

Set-TextValue $ws.Range('D26') '9.163'
Set-TextValue $ws.Range('E26') '  -1.65%  '

Set-TextValue $ws.Range('E27') '  +0.24%  '

Set-TextValue $ws.Range('D28') '1.916'
Set-TextValue $ws.Range('E28') '  -0.13%  '

Set-TextValue $ws.Range('E29') '  +3.21%  '

Set-TextValue $ws.Range('D30') '0.09680'
Set-TextValue $ws.Range('E30') '  +1.06%  '

Set-TextValue $ws.Range('D31') '4.359'
Set-TextValue $ws.Range('E31') '  -1.10%  '

Set-TextValue $ws.Range('D32') '1.476'
Set-TextValue $ws.Range('E32') '  +0.47%  '

Set-TextValue $ws.Range('E33') '  -1.84%  '

Set-TextValue $ws.Range('D34') '0.04701'
Set-TextValue $ws.Range('E34') '  +1.01%  '

Set-TextValue $ws.Range('E35') '  +2.41%  '

Set-TextValue $ws.Range('D36') '0.7036'
Set-TextValue $ws.Range('E36') '  +0.56%  '

Set-TextValue $ws.Range('D37') '2.708'
Set-TextValue $ws.Range('E37') '  -0.01%  '

Set-TextValue $ws.Range('E38') '  +0.07%  '

Set-TextValue $ws.Range('D39') '2.610'
Set-TextValue $ws.Range('E39') '  +3.89%  '

Set-TextValue $ws.Range('D40') '6.290'
Set-TextValue $ws.Range('E40') '  -1.67%  '

Set-TextValue $ws.Range('D41') '74.17'
Set-TextValue $ws.Range('E41') '  +2.87%  '

Set-TextValue $ws.Range('D42') '1.943'
Set-TextValue $ws.Range('E42') '  +0.64%  '

Set-TextValue $ws.Range('D43') '0.8483'
Set-TextValue $ws.Range('E43') '  -0.71%  '

Set-TextValue $ws.Range('D44') '1.000'
Set-TextValue $ws.Range('E44') '  -0.01%  '

Set-TextValue $ws.Range('D45') '0.4152'
Set-TextValue $ws.Range('E45') '  -0.35%  '

Set-TextValue $ws.Range('D46') '103.16'
Set-TextValue $ws.Range('E46') '  -0.63%  '

Set-TextValue $ws.Range('D47') '981.86'
Set-TextValue $ws.Range('E47') '  -1.71%  '

Set-TextValue $ws.Range('D48') '7.135'
Set-TextValue $ws.Range('E48') '  -0.76%  '

Set-TextValue $ws.Range('D49') '9.290'
Set-TextValue $ws.Range('E49') '  +2.31%  '

Set-TextValue $ws.Range('D50') '33.92'
Set-TextValue $ws.Range('E50') '  -0.10%  '

Set-TextValue $ws.Range('D51') '0.05644'
Set-TextValue $ws.Range('E51') '  +0.14%  '
